$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 4999.857
$ws.Range("J3").Value = 4999.857
$ws.Range("L3").Value = 4999.857
$ws.Range("N3").Value = -5227.857
$ws.Range("H40").Value = 1550.125
$ws.Range("I40").Value = 933.6667
$ws.Range("J40").Value = 1920
$ws.Range("K40").Value = 933.6667
$ws.Range("L40").Value = 1920
$ws.Range("M40").Value = -758.6667
$ws.Range("N40").Value = -2270
$ws.Range("H51").Value = 3144.5557
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 3287.625
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 3287.625
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -4255.625
$ws.Range("H54").Value = 3450
$ws.Range("I54").Value = 2000
$ws.Range("J54").Value = 4900
$ws.Range("K54").Value = 2000
$ws.Range("L54").Value = 4900
$ws.Range("M54").Value = -1514
$ws.Range("N54").Value = -5872
$ws.Range("H62").Value = 3283.3333
$ws.Range("I62").Value = 2940
$ws.Range("K62").Value = 2940
$ws.Range("M62").Value = -2316
$ws.Range("H65").Value = 3283.3333
$ws.Range("I65").Value = 2940
$ws.Range("K65").Value = 14700
$ws.Range("M65").Value = -11580
$ws.Range("H86").Value = 3475.682
$ws.Range("I86").Value = 2115.5881
$ws.Range("J86").Value = 8100
$ws.Range("K86").Value = 2115.5881
$ws.Range("L86").Value = 8100
$ws.Range("M86").Value = -992.5880999999999
$ws.Range("N86").Value = -10346
$ws.Range("H89").Value = 3475.682
$ws.Range("I89").Value = 2115.5881
$ws.Range("J89").Value = 8100
$ws.Range("K89").Value = 10577.9405
$ws.Range("L89").Value = 40500
$ws.Range("M89").Value = -4961.940500000001
$ws.Range("N89").Value = -51732
$ws.Range("H102").Value = 4999.857
$ws.Range("J102").Value = 4999.857
$ws.Range("L102").Value = 4999.857
$ws.Range("N102").Value = -11489.857
$ws.Range("H116").Value = 2755.074
$ws.Range("I116").Value = 2327.625
$ws.Range("J116").Value = 3376.818
$ws.Range("K116").Value = 2327.625
$ws.Range("L116").Value = 3376.818
$ws.Range("M116").Value = 1114.375
$ws.Range("N116").Value = -10260.818
$ws.Range("H132").Value = 8550214
$ws.Range("I132").Value = 11908278
$ws.Range("J132").Value = 2415.6365
$ws.Range("K132").Value = 35724834
$ws.Range("L132").Value = 7246.9095
$ws.Range("M132").Value = -35722304
$ws.Range("N132").Value = -12306.9095
$ws.Range("H137").Value = 1697.5416
$ws.Range("I137").Value = 1341.1666
$ws.Range("K137").Value = 4023.4998
$ws.Range("M137").Value = -1473.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3419.68
$ws.Range("I32").Value = 3659.182
$ws.Range("K32").Value = 3659.182
$ws.Range("M32").Value = -3372.182
$ws.Range("H45").Value = 2185.125
$ws.Range("I45").Value = 2104.4285
$ws.Range("K45").Value = 2104.4285
$ws.Range("M45").Value = -1727.4285
$ws.Range("H132").Value = 2429.8262
$ws.Range("I132").Value = 2155.2974
$ws.Range("J132").Value = 3558.4443
$ws.Range("K132").Value = 6465.8922
$ws.Range("L132").Value = 10675.3329
$ws.Range("M132").Value = -3935.8922
$ws.Range("N132").Value = -15735.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8929238
$ws.Range("I94").Value = 12500519
$ws.Range("J94").Value = 1036.125
$ws.Range("K94").Value = 12500519
$ws.Range("L94").Value = 1036.125
$ws.Range("M94").Value = -12500068
$ws.Range("N94").Value = -1938.125
$ws.Range("H134").Value = 5330.0835
$ws.Range("I134").Value = 1233
$ws.Range("J134").Value = 12158.556
$ws.Range("K134").Value = 3699
$ws.Range("L134").Value = 36475.66800000001
$ws.Range("M134").Value = -1164
$ws.Range("N134").Value = -41545.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1024
$ws.Range("H31").Value = 994.3019
$ws.Range("I31").Value = 770.7931
$ws.Range("J31").Value = 1264.375
$ws.Range("K31").Value = 770.7931
$ws.Range("L31").Value = 1264.375
$ws.Range("M31").Value = -475.7931
$ws.Range("N31").Value = -1854.375
$ws.Range("H34").Value = 994.3019
$ws.Range("I34").Value = 770.7931
$ws.Range("J34").Value = 1264.375
$ws.Range("K34").Value = 770.7931
$ws.Range("L34").Value = 1264.375
$ws.Range("M34").Value = -568.7931
$ws.Range("N34").Value = -1668.375
$ws.Range("H132").Value = 8241.3125
$ws.Range("I132").Value = 8989.462
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 26968.386
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -24438.386
$ws.Range("N132").Value = -20058.0005
$ws.Range("H134").Value = 9805156
$ws.Range("I134").Value = 11905903
$ws.Range("K134").Value = 35717709
$ws.Range("M134").Value = -35715174

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3500
$ws.Range("J55").Value = 3500
$ws.Range("L55").Value = 10500
$ws.Range("N55").Value = -10854
$ws.Range("H68").Value = 2189.5122
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 2276.1538
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 6828.4614
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -8450.4614
$ws.Range("H69").Value = 2440.9333
$ws.Range("J69").Value = 2401
$ws.Range("L69").Value = 7203
$ws.Range("N69").Value = -8825
$ws.Range("H71").Value = 2189.5122
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 2276.1538
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 20485.3842
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -28597.3842
$ws.Range("H72").Value = 2440.9333
$ws.Range("J72").Value = 2401
$ws.Range("L72").Value = 21609
$ws.Range("N72").Value = -29721
$ws.Range("H107").Value = 7057.8125
$ws.Range("I107").Value = 621.8570999999999
$ws.Range("J107").Value = 12063.556
$ws.Range("K107").Value = 1865.5713
$ws.Range("L107").Value = 36190.66800000001
$ws.Range("M107").Value = 54.42870000000016
$ws.Range("N107").Value = -40030.66800000001
$ws.Range("H131").Value = 19609290
$ws.Range("J131").Value = 1594.186
$ws.Range("L131").Value = 4782.558
$ws.Range("N131").Value = -14862.558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23689030
$ws.Range("I70").Value = 27781618
$ws.Range("J70").Value = 20005702
$ws.Range("K70").Value = 27781618
$ws.Range("L70").Value = 20005702
$ws.Range("M70").Value = -27781348
$ws.Range("N70").Value = -20006242
$ws.Range("H73").Value = 23689030
$ws.Range("I73").Value = 27781618
$ws.Range("J73").Value = 20005702
$ws.Range("K73").Value = 27781618
$ws.Range("L73").Value = 20005702
$ws.Range("M73").Value = -27780682
$ws.Range("N73").Value = -20007574
$ws.Range("H97").Value = 1947.6923
$ws.Range("I97").Value = 2402.111
$ws.Range("J97").Value = 925.25
$ws.Range("K97").Value = 2402.111
$ws.Range("L97").Value = 925.25
$ws.Range("M97").Value = -1906.111
$ws.Range("N97").Value = -1917.25
$ws.Range("H98").Value = 15257
$ws.Range("J98").Value = 15257
$ws.Range("L98").Value = 15257
$ws.Range("N98").Value = -21247
$ws.Range("H132").Value = 2213.4443
$ws.Range("I132").Value = 1742.4
$ws.Range("J132").Value = 4568.6665
$ws.Range("K132").Value = 5227.200000000001
$ws.Range("L132").Value = 13705.9995
$ws.Range("M132").Value = -2697.200000000001
$ws.Range("N132").Value = -18765.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2072.2222
$ws.Range("I16").Value = 2018.75
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2018.75
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1848.75
$ws.Range("N16").Value = -2840
$ws.Range("H22").Value = 796
$ws.Range("I22").Value = 415.25
$ws.Range("J22").Value = 1049.8334
$ws.Range("K22").Value = 415.25
$ws.Range("L22").Value = 1049.8334
$ws.Range("M22").Value = -120.25
$ws.Range("N22").Value = -1639.8334
$ws.Range("H27").Value = 796
$ws.Range("I27").Value = 415.25
$ws.Range("J27").Value = 1049.8334
$ws.Range("K27").Value = 415.25
$ws.Range("L27").Value = 1049.8334
$ws.Range("M27").Value = -308.25
$ws.Range("N27").Value = -1263.8334
$ws.Range("H74").Value = 35000
$ws.Range("J74").Value = 35000
$ws.Range("L74").Value = 35000
$ws.Range("N74").Value = -36996
$ws.Range("H77").Value = 35000
$ws.Range("J77").Value = 35000
$ws.Range("L77").Value = 105000
$ws.Range("N77").Value = -114984
$ws.Range("H100").Value = 2042
$ws.Range("I100").Value = 1966.5
$ws.Range("J100").Value = 2117.5
$ws.Range("K100").Value = 1966.5
$ws.Range("L100").Value = 2117.5
$ws.Range("M100").Value = -1425.5
$ws.Range("N100").Value = -3199.5
$ws.Range("H106").Value = 29999.5
$ws.Range("J106").Value = 29999.5
$ws.Range("L106").Value = 29999.5
$ws.Range("N106").Value = -32523.5
$ws.Range("H122").Value = 94449450
$ws.Range("I122").Value = 94449450
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 283348350
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -283345900
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 33107.5
$ws.Range("J135").Value = 33107.5
$ws.Range("L135").Value = 33107.5
$ws.Range("N135").Value = -43247.5
$ws.Range("H136").Value = 8232.200000000001
$ws.Range("I136").Value = 15583.429
$ws.Range("J136").Value = 1799.875
$ws.Range("K136").Value = 46750.287
$ws.Range("L136").Value = 5399.625
$ws.Range("M136").Value = -44200.287
$ws.Range("N136").Value = -10499.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33341446
$ws.Range("I62").Value = 50005740
$ws.Range("J62").Value = 12859.8
$ws.Range("K62").Value = 50005740
$ws.Range("L62").Value = 12859.8
$ws.Range("M62").Value = -50005116
$ws.Range("N62").Value = -14107.8
$ws.Range("H65").Value = 33341446
$ws.Range("I65").Value = 50005740
$ws.Range("J65").Value = 12859.8
$ws.Range("K65").Value = 250028700
$ws.Range("L65").Value = 64299
$ws.Range("M65").Value = -250025580
$ws.Range("N65").Value = -70539
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 2632
$ws.Range("I132").Value = 2256.6333
$ws.Range("J132").Value = 4240.7144
$ws.Range("K132").Value = 6769.8999
$ws.Range("L132").Value = 12722.1432
$ws.Range("M132").Value = -4239.8999
$ws.Range("N132").Value = -17782.1432
$ws.Range("H136").Value = 622.0476
$ws.Range("I136").Value = 400.92856
$ws.Range("J136").Value = 1064.2858
$ws.Range("K136").Value = 1202.78568
$ws.Range("L136").Value = 3192.8574
$ws.Range("M136").Value = 1347.21432
$ws.Range("N136").Value = -8292.857400000001
